# Weekly update: a new price-record row is inserted at the top of the
# data (row 5), pushing all existing price rows down by one. The new
# row captures the newest week's observation; everything below keeps
# its previous values (just shifted down one row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 5; this shifts rows 5:40 down to 6:41
# and extends the used range to A1:R41.
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row 5 with this week's record.
$ws.Cells.Item(5, 1).Value = 1
$ws.Cells.Item(5, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(5, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(5, 4).Value = 45051
$ws.Cells.Item(5, 5).Value = 15
$ws.Cells.Item(5, 6).Value = 100112003
$ws.Cells.Item(5, 7).Value = "Ajo"
$ws.Cells.Item(5, 8).Value = "Chino"
$ws.Cells.Item(5, 9).Value = "Primera"
$ws.Cells.Item(5, 10).Value = 750
$ws.Cells.Item(5, 11).Value = 16000
$ws.Cells.Item(5, 12).Value = 17000
$ws.Cells.Item(5, 13).Value = 16333
$ws.Cells.Item(5, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(5, 15).Value = "China"
$ws.Cells.Item(5, 16).Value = 1633
$ws.Cells.Item(5, 17).Value = 10
$ws.Cells.Item(5, 18).Value = "Hortaliza"
